$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell data for rows 2-8 (post-edit state).
#   T="t"    => plain text value (string)
#   T="tnum" => text value that looks like a number (needs a leading
#               apostrophe so Excel keeps storing it as text, not a number)
#   T="v"    => numeric value
$cellData = @(
    @{R=2; C=1; T="t"; V="Jersey"}
    @{R=2; C=2; T="tnum"; V="6"}
    @{R=2; C=3; T="t"; V="Investments & Asset Management"}
    @{R=2; C=4; T="v"; V="0.135"}
    @{R=2; C=5; T="v"; V="-0.104"}
    @{R=2; C=6; T="v"; V="-0.042"}
    @{R=2; C=7; T="v"; V="0.35252771129608"}
    @{R=2; C=8; T="v"; V="0.3524484778187105"}
    @{R=2; C=9; T="v"; V="0.2301936727575701"}
    @{R=2; C=10; T="v"; V="0.216707348124562"}
    @{R=2; C=11; T="v"; V="233.189"}
    @{R=2; C=12; T="v"; V="0.1904780964361272"}
    @{R=2; C=13; T="v"; V="327.051"}
    @{R=2; C=14; T="v"; V="0.07449230138483966"}
    @{R=2; C=15; T="v"; V="1.402514698377711"}
    @{R=2; C=16; T="v"; V="180.737"}
    @{R=2; C=17; T="v"; V="0.0411664085276968"}
    @{R=2; C=18; T="v"; V="0.7750665768968519"}
    @{R=2; C=19; T="v"; V="146.314"}
    @{R=2; C=20; T="v"; V="0.4473736512042464"}
    @{R=2; C=21; T="v"; V="328.818"}
    @{R=2; C=22; T="v"; V="0.07489477040816327"}
    @{R=2; C=23; T="v"; V="-0.007932115875592591"}
    @{R=2; C=24; T="v"; V="0.03802262921432312"}
    @{R=2; C=25; T="v"; V="-0.04595474508991572"}
    @{R=2; C=26; T="v"; V="0.7857782878386428"}
    @{R=2; C=27; T="v"; V="0.01021456981369061"}
    @{R=2; C=28; T="v"; V="0.03643585029976566"}
    @{R=2; C=29; T="v"; V="-0.02612323771483396"}
    @{R=2; C=30; T="v"; V="521.087"}
    @{R=2; C=31; T="v"; V="0"}
    @{R=2; C=32; T="v"; V="521.087"}
    @{R=2; C=33; T="v"; V="192.269"}
    @{R=2; C=34; T="v"; V="0.1060955673913013"}
    @{R=2; C=35; T="v"; V="0.1988512059018037"}
    @{R=2; C=36; T="v"; V="0.04195568128529466"}
    @{R=2; C=37; T="v"; V="0.08389911457544698"}
    @{R=2; C=38; T="v"; V="9.108000000000001"}
    @{R=2; C=39; T="v"; V="2.617000000000001"}
    @{R=2; C=40; T="v"; V="1.296952545933326"}
    @{R=2; C=41; T="v"; V="30.9409310496267"}
    @{R=2; C=42; T="v"; V="0.4785453658488021"}
    @{R=2; C=43; T="v"; V="107.6843714176538"}
    @{R=3; C=1; T="t"; V="Jersey"}
    @{R=3; C=2; T="t"; V="Sanne Group plc (LSE:SNN)"}
    @{R=3; C=3; T="t"; V="Investments & Asset Management"}
    @{R=3; C=4; T="v"; V="0.333"}
    @{R=3; C=7; T="v"; V="0.2824687199230029"}
    @{R=3; C=8; T="v"; V="0.2820019249278152"}
    @{R=3; C=9; T="v"; V="0.1732435033686237"}
    @{R=3; C=10; T="v"; V="0.1247192813602823"}
    @{R=3; C=11; T="v"; V="18.9"}
    @{R=3; C=12; T="v"; V="0.09095283926852742"}
    @{R=3; C=13; T="v"; V="25.414"}
    @{R=3; C=14; T="v"; V="0.02075967979088384"}
    @{R=3; C=15; T="v"; V="1.344656084656085"}
    @{R=3; C=16; T="v"; V="25.3"}
    @{R=3; C=17; T="v"; V="0.02066655775200131"}
    @{R=3; C=18; T="v"; V="1.338624338624339"}
    @{R=3; C=19; T="v"; V="0.1140000000000008"}
    @{R=3; C=20; T="v"; V="0.004485716534193781"}
    @{R=3; C=21; T="v"; V="63.8"}
    @{R=3; C=22; T="v"; V="0.05211566737461198"}
    @{R=3; C=23; T="v"; V="0.08080376229157758"}
    @{R=3; C=24; T="v"; V="0.03992933431246774"}
    @{R=3; C=25; T="v"; V="0.04087442797910985"}
    @{R=3; C=26; T="v"; V="1.821209465381244"}
    @{R=3; C=27; T="v"; V="0.2271399357288927"}
    @{R=3; C=28; T="v"; V="0.03675940427990798"}
    @{R=3; C=29; T="v"; V="0.1903805314489848"}
    @{R=3; C=30; T="v"; V="221"}
    @{R=3; C=31; T="v"; V="0"}
    @{R=3; C=32; T="v"; V="221"}
    @{R=3; C=33; T="v"; V="157.2"}
    @{R=3; C=34; T="v"; V="0.152920011071132"}
    @{R=3; C=35; T="v"; V="0.4933035714285715"}
    @{R=3; C=36; T="v"; V="0.1137975966410887"}
    @{R=3; C=37; T="v"; V="0.4091618948464341"}
    @{R=3; C=38; T="v"; V="5.66"}
    @{R=3; C=39; T="v"; V="5.472"}
    @{R=3; C=40; T="v"; V="3.689482470784641"}
    @{R=3; C=41; T="v"; V="6.360424028268551"}
    @{R=3; C=42; T="v"; V="2.624373956594324"}
    @{R=3; C=43; T="v"; V="6.578947368421052"}
    @{R=4; C=1; T="t"; V="Jersey"}
    @{R=4; C=2; T="t"; V="Man Group plc (LSE:EMG)"}
    @{R=4; C=3; T="t"; V="Investments & Asset Management"}
    @{R=4; C=4; T="v"; V="-0.0482"}
    @{R=4; C=5; T="v"; V="-0.104"}
    @{R=4; C=6; T="v"; V="-0.042"}
    @{R=4; C=7; T="v"; V="0.3734939759036144"}
    @{R=4; C=8; T="v"; V="0.3734939759036144"}
    @{R=4; C=9; T="v"; V="0.2570281124497992"}
    @{R=4; C=10; T="v"; V="0.2386689615605278"}
    @{R=4; C=11; T="v"; V="234"}
    @{R=4; C=12; T="v"; V="0.2349397590361446"}
    @{R=4; C=13; T="v"; V="287"}
    @{R=4; C=14; T="v"; V="0.1067430356677948"}
    @{R=4; C=15; T="v"; V="1.226495726495727"}
    @{R=4; C=16; T="v"; V="147"}
    @{R=4; C=17; T="v"; V="0.05467326217130956"}
    @{R=4; C=18; T="v"; V="0.6282051282051282"}
    @{R=4; C=19; T="v"; V="140"}
    @{R=4; C=20; T="v"; V="0.4878048780487805"}
    @{R=4; C=21; T="v"; V="236"}
    @{R=4; C=22; T="v"; V="0.08777476103693235"}
    @{R=4; C=23; T="v"; V="0.1482889733840304"}
    @{R=4; C=24; T="v"; V="0.03773515214231193"}
    @{R=4; C=25; T="v"; V="0.1105538212417185"}
    @{R=4; C=26; T="v"; V="0.848381601362862"}
    @{R=4; C=27; T="v"; V="0.202482355804332"}
    @{R=4; C=28; T="v"; V="0.03611229631962334"}
    @{R=4; C=29; T="v"; V="0.1663700594847086"}
    @{R=4; C=30; T="v"; V="258"}
    @{R=4; C=31; T="v"; V="0"}
    @{R=4; C=32; T="v"; V="258"}
    @{R=4; C=33; T="v"; V="22"}
    @{R=4; C=34; T="v"; V="0.0875555706383412"}
    @{R=4; C=35; T="v"; V="0.1405228758169935"}
    @{R=4; C=36; T="v"; V="0.008115984800973919"}
    @{R=4; C=37; T="v"; V="0.01375"}
    @{R=4; C=38; T="v"; V="1"}
    @{R=4; C=39; T="v"; V="-5"}
    @{R=4; C=40; T="v"; V="0.7565982404692082"}
    @{R=4; C=41; T="v"; V="256"}
    @{R=4; C=42; T="v"; V="0.06451612903225806"}
    @{R=4; C=43; T="v"; V="-51.2"}
    @{R=5; C=1; T="t"; V="Jersey"}
    @{R=5; C=2; T="t"; V="MJ Hudson Group plc (AIM:MJH)"}
    @{R=5; C=3; T="t"; V="Investments & Asset Management"}
    @{R=5; C=7; T="v"; V="0.03181159420289855"}
    @{R=5; C=8; T="v"; V="0.03181159420289855"}
    @{R=5; C=9; T="v"; V="0.001956521739130435"}
    @{R=5; C=10; T="v"; V="0.001956521739130435"}
    @{R=5; C=11; T="v"; V="-8.92"}
    @{R=5; C=12; T="v"; V="-0.3231884057971015"}
    @{R=5; C=13; T="v"; V="-0"}
    @{R=5; C=14; T="v"; V="-0"}
    @{R=5; C=15; T="v"; V="0"}
    @{R=5; C=16; T="v"; V="-0"}
    @{R=5; C=17; T="v"; V="-0"}
    @{R=5; C=18; T="v"; V="0"}
    @{R=5; C=19; T="v"; V="0"}
    @{R=5; C=21; T="v"; V="16.6"}
    @{R=5; C=22; T="v"; V="0.1474245115452931"}
    @{R=5; C=23; T="v"; V="-0.9019211324570272"}
    @{R=5; C=24; T="v"; V="0.03831010628633431"}
    @{R=5; C=25; T="v"; V="-0.9402312387433616"}
    @{R=5; C=26; T="v"; V="10.82352941176472"}
    @{R=5; C=27; T="v"; V="0.02117647058823532"}
    @{R=5; C=28; T="v"; V="0.03743011888674538"}
    @{R=5; C=29; T="v"; V="-0.01625364829851007"}
    @{R=5; C=30; T="v"; V="13.3"}
    @{R=5; C=31; T="v"; V="0"}
    @{R=5; C=32; T="v"; V="13.3"}
    @{R=5; C=33; T="v"; V="-3.300000000000001"}
    @{R=5; C=34; T="v"; V="0.1056393963463066"}
    @{R=5; C=35; T="v"; V="0.2074882995319813"}
    @{R=5; C=36; T="v"; V="-0.030192131747484"}
    @{R=5; C=37; T="v"; V="-0.06947368421052633"}
    @{R=5; C=38; T="v"; V="1.15"}
    @{R=5; C=39; T="v"; V="1.15"}
    @{R=5; C=40; T="v"; V="15.14806378132119"}
    @{R=5; C=41; T="v"; V="0.04695652173913044"}
    @{R=5; C=42; T="v"; V="-3.758542141230069"}
    @{R=5; C=43; T="v"; V="0.04695652173913044"}
    @{R=6; C=1; T="t"; V="Jersey"}
    @{R=6; C=2; T="t"; V="TMT Investments PLC (AIM:TMT)"}
    @{R=6; C=3; T="t"; V="Investments & Asset Management"}
    @{R=6; C=6; T="v"; V="0.135"}
    @{R=6; C=7; T="v"; V="0"}
    @{R=6; C=8; T="v"; V="0"}
    @{R=6; C=9; T="v"; V="-0.05779816513761468"}
    @{R=6; C=10; T="v"; V="-0.05779816513761468"}
    @{R=6; C=11; T="v"; V="0.246"}
    @{R=6; C=12; T="v"; V="0.2256880733944954"}
    @{R=6; C=13; T="v"; V="0.037"}
    @{R=6; C=14; T="v"; V="0.0001589347079037801"}
    @{R=6; C=15; T="v"; V="0.1504065040650406"}
    @{R=6; C=16; T="v"; V="0.037"}
    @{R=6; C=17; T="v"; V="0.0001589347079037801"}
    @{R=6; C=18; T="v"; V="0.1504065040650406"}
    @{R=6; C=19; T="v"; V="0"}
    @{R=6; C=20; T="v"; V="0"}
    @{R=6; C=21; T="v"; V="9.19"}
    @{R=6; C=22; T="v"; V="0.03947594501718213"}
    @{R=6; C=23; T="v"; V="0.002305529522024367"}
    @{R=6; C=24; T="v"; V="0.03524549617030377"}
    @{R=6; C=25; T="v"; V="-0.0329399666482794"}
    @{R=6; C=26; T="v"; V="0.0129300118623962"}
    @{R=6; C=27; T="v"; V="-0.0007473309608540924"}
    @{R=6; C=28; T="v"; V="0.03524549617030377"}
    @{R=6; C=29; T="v"; V="-0.03599282713115786"}
    @{R=6; C=30; T="v"; V="0"}
    @{R=6; C=31; T="v"; V="0"}
    @{R=6; C=32; T="v"; V="0"}
    @{R=6; C=33; T="v"; V="-9.19"}
    @{R=6; C=34; T="v"; V="0"}
    @{R=6; C=35; T="v"; V="0"}
    @{R=6; C=36; T="v"; V="-0.04109834086132105"}
    @{R=6; C=37; T="v"; V="-0.09998911979109999"}
    @{R=6; C=38; T="v"; V="0"}
    @{R=6; C=42; T="v"; V="-0.303"}
    @{R=6; C=43; T="v"; V="0.2079207920792079"}
    @{R=7; C=1; T="t"; V="Jersey"}
    @{R=7; C=2; T="t"; V="Westmount Energy Limited (AIM:WTE)"}
    @{R=7; C=3; T="t"; V="Investments & Asset Management"}
    @{R=7; C=12; T="v"; V="-0.137"}
    @{R=7; C=13; T="v"; V="-0"}
    @{R=7; C=14; T="v"; V="-0"}
    @{R=7; C=15; T="v"; V="0"}
    @{R=7; C=16; T="v"; V="-0"}
    @{R=7; C=17; T="v"; V="-0"}
    @{R=7; C=18; T="v"; V="0"}
    @{R=7; C=20; T="v"; V="0"}
    @{R=7; C=21; T="v"; V="3.02"}
    @{R=7; C=22; T="v"; V="0.08779069767441861"}
    @{R=7; C=23; T="v"; V="-0.01816976127320955"}
    @{R=7; C=24; T="v"; V="0.03561280595620312"}
    @{R=7; C=25; T="v"; V="-0.05378256722941267"}
    @{R=7; C=26; T="v"; V="0"}
    @{R=7; C=27; T="v"; V="-0.04477972238986119"}
    @{R=7; C=28; T="v"; V="0.03547582551934495"}
    @{R=7; C=29; T="v"; V="-0.08025554790920614"}
    @{R=7; C=30; T="v"; V="0.487"}
    @{R=7; C=31; T="v"; V="0"}
    @{R=7; C=32; T="v"; V="0.487"}
    @{R=7; C=33; T="v"; V="-2.533"}
    @{R=7; C=34; T="v"; V="0.01395935448734486"}
    @{R=7; C=35; T="v"; V="0.02737954686006634"}
    @{R=7; C=36; T="v"; V="-0.07948661624878402"}
    @{R=7; C=37; T="v"; V="-0.1715311166790817"}
    @{R=7; C=38; T="v"; V="0.068"}
    @{R=7; C=39; T="v"; V="0.068"}
    @{R=7; C=41; T="v"; V="-5.455882352941176"}
    @{R=7; C=43; T="v"; V="-5.455882352941176"}
    @{R=8; C=1; T="t"; V="Jersey"}
    @{R=8; C=2; T="t"; V="EJF Investments Limited (LSE:EJFI)"}
    @{R=8; C=3; T="t"; V="Investments & Asset Management"}
    @{R=8; C=7; T="v"; V="-0"}
    @{R=8; C=8; T="v"; V="-0"}
    @{R=8; C=9; T="v"; V="1.187651331719128"}
    @{R=8; C=10; T="v"; V="1.187651331719128"}
    @{R=8; C=11; T="v"; V="-10.9"}
    @{R=8; C=12; T="v"; V="1.319612590799031"}
    @{R=8; C=13; T="v"; V="14.6"}
    @{R=8; C=14; T="v"; V="0.1494370522006141"}
    @{R=8; C=15; T="v"; V="-1.339449541284404"}
    @{R=8; C=16; T="v"; V="8.4"}
    @{R=8; C=17; T="v"; V="0.08597748208802457"}
    @{R=8; C=18; T="v"; V="-0.7706422018348624"}
    @{R=8; C=19; T="v"; V="6.200000000000001"}
    @{R=8; C=20; T="v"; V="0.4246575342465754"}
    @{R=8; C=21; T="v"; V="0.208"}
    @{R=8; C=22; T="v"; V="0.002128966223132037"}
    @{R=8; C=23; T="v"; V="-0.070458952811894"}
    @{R=8; C=24; T="v"; V="0.04276092648370802"}
    @{R=8; C=25; T="v"; V="-0.113219879295602"}
    @{R=8; C=26; T="v"; V="-0.04726779552386566"}
    @{R=8; C=27; T="v"; V="-0.05613766030134651"}
    @{R=8; C=28; T="v"; V="0.03989028982109741"}
    @{R=8; C=29; T="v"; V="-0.09602795012244392"}
    @{R=8; C=30; T="v"; V="28.3"}
    @{R=8; C=31; T="v"; V="0"}
    @{R=8; C=32; T="v"; V="28.3"}
    @{R=8; C=33; T="v"; V="28.092"}
    @{R=8; C=34; T="v"; V="0.2246031746031746"}
    @{R=8; C=35; T="v"; V="0.1843648208469055"}
    @{R=8; C=36; T="v"; V="0.2233210379038413"}
    @{R=8; C=37; T="v"; V="0.1832580956605694"}
    @{R=8; C=38; T="v"; V="1.23"}
    @{R=8; C=39; T="v"; V="1.23"}
    @{R=8; C=41; T="v"; V="-7.975609756097562"}
    @{R=8; C=43; T="v"; V="-7.975609756097562"}
)

# Clear out all existing data in rows 2:9, columns A:AQ, so stale values
# from the old row layout don't linger after the row shuffle/removal.
$clearRange = $ws.Range($ws.Cells.Item(2,1), $ws.Cells.Item(9,43))
$clearRange.ClearContents()

foreach ($cell in $cellData) {
    $target = $ws.Cells.Item($cell.R, $cell.C)
    if ($cell.T -eq "tnum") {
        $target.Value = "'" + $cell.V
    } elseif ($cell.T -eq "t") {
        $target.Value = $cell.V
    } else {
        $target.Value = [double]$cell.V
    }
}

# Row 9 (old "MJ Hudson" row before the reshuffle) no longer exists in the
# updated dataset; delete it so the sheet dimension shrinks to A1:AQ8.
$ws.Rows.Item(9).Delete()
